# Update excess mortality, prepare data for nowcast monkeypox
#
# Applies the edits described by the diff:
#  - a handful of single-cell data corrections (+1 adjustments) scattered
#    through rows 108-135, whose dependent "% change" formulas recalc
#    automatically,
#  - a bigger set of data corrections on rows 133-135 (new provisional
#    week numbers revised),
#  - a brand-new week of data (row 136, "2022 week 30") with its 12
#    "% change vs baseline" formulas,
#  - and the sheet view (scroll position / selection) is moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Small, isolated data corrections (dependent formula cells recalc
#    automatically once the underlying numbers change).
# ---------------------------------------------------------------------
$ws.Range("W108").Value  = 467
$ws.Range("W109").Value  = 515
$ws.Range("R116").Value  = 102
$ws.Range("W123").Value  = 447
$ws.Range("W124").Value  = 485
$ws.Range("V126").Value  = 197
$ws.Range("W127").Value  = 431
$ws.Range("Z127").Value  = 447
$ws.Range("W128").Value  = 456
$ws.Range("S130").Value  = 203
$ws.Range("S131").Value  = 222
$ws.Range("W131").Value  = 444
$ws.Range("X131").Value  = 700
$ws.Range("AA132").Value = 258

# ---------------------------------------------------------------------
# 2. Row 133 revisions.
# ---------------------------------------------------------------------
$ws.Range("P133").Value  = 96
$ws.Range("S133").Value  = 184
$ws.Range("U133").Value  = 382
$ws.Range("V133").Value  = 195
$ws.Range("X133").Value  = 601
$ws.Range("AA133").Value = 208

# ---------------------------------------------------------------------
# 3. Row 134 revisions.
# ---------------------------------------------------------------------
$ws.Range("U134").Value  = 369
$ws.Range("V134").Value  = 216
$ws.Range("W134").Value  = 505
$ws.Range("X134").Value  = 619
$ws.Range("Y134").Value  = 78
$ws.Range("Z134").Value  = 444
$ws.Range("AA134").Value = 243

# ---------------------------------------------------------------------
# 4. Row 135 revisions.
# ---------------------------------------------------------------------
$ws.Range("P135").Value  = 122
$ws.Range("Q135").Value  = 128
$ws.Range("R135").Value  = 128
$ws.Range("S135").Value  = 202
$ws.Range("T135").Value  = 58
$ws.Range("U135").Value  = 400
$ws.Range("V135").Value  = 193
$ws.Range("W135").Value  = 498
$ws.Range("X135").Value  = 625
$ws.Range("Y135").Value  = 96
$ws.Range("Z135").Value  = 451
$ws.Range("AA135").Value = 243

# ---------------------------------------------------------------------
# 5. Row 136: newly added week ("2022 week 30"), provinces P:AA plus the
#    year/week labels and the twelve "% change vs baseline" formulas
#    (same pattern used by every other row: ROUND((new-base)/base*100,2)).
# ---------------------------------------------------------------------
$ws.Range("P136").Value  = 116
$ws.Range("Q136").Value  = 135
$ws.Range("R136").Value  = 98
$ws.Range("S136").Value  = 228
$ws.Range("T136").Value  = 58
$ws.Range("U136").Value  = 421
$ws.Range("V136").Value  = 201
$ws.Range("W136").Value  = 539
$ws.Range("X136").Value  = 637
$ws.Range("Y136").Value  = 57
$ws.Range("Z136").Value  = 475
$ws.Range("AA136").Value = 248
$ws.Range("AC136").Value = 2022
$ws.Range("AD136").Value = 30

$ws.Range("AE136").Formula = "=ROUND((P136-B136)/B136*100,2)"
$ws.Range("AF136").Formula = "=ROUND((Q136-C136)/C136*100,2)"
$ws.Range("AG136").Formula = "=ROUND((R136-D136)/D136*100,2)"
$ws.Range("AH136").Formula = "=ROUND((S136-E136)/E136*100,2)"
$ws.Range("AI136").Formula = "=ROUND((T136-F136)/F136*100,2)"
$ws.Range("AJ136").Formula = "=ROUND((U136-G136)/G136*100,2)"
$ws.Range("AK136").Formula = "=ROUND((V136-H136)/H136*100,2)"
$ws.Range("AL136").Formula = "=ROUND((W136-I136)/I136*100,2)"
$ws.Range("AM136").Formula = "=ROUND((X136-J136)/J136*100,2)"
$ws.Range("AN136").Formula = "=ROUND((Y136-K136)/K136*100,2)"
$ws.Range("AO136").Formula = "=ROUND((Z136-L136)/L136*100,2)"
$ws.Range("AP136").Formula = "=ROUND((AA136-M136)/M136*100,2)"

# ---------------------------------------------------------------------
# 6. Sheet view: scroll position + active selection moved.
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 103
$win.ScrollColumn = 5
$ws.Range("AI105").Select() | Out-Null
